$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 50 (pushes existing rows 50..158 down to 51..159)
$ws.Rows.Item(50).Insert()

# Populate the new row's Key / Comment / English columns
$ws.Range("B50").Value = "strFileHeader29"
$ws.Range("C50").Value = "Field description in exported file"
$ws.Range("D50").Value = "Differentiation algorithm"

# Style B50: left/center alignment, no wrap (matches style index 2)
$ws.Range("B50").HorizontalAlignment = -4131
$ws.Range("B50").VerticalAlignment = -4108
$ws.Range("B50").WrapText = $false

# Style C50:D50: left/center alignment, no wrap, unlocked (matches style index 3)
$ws.Range("C50:D50").HorizontalAlignment = -4131
$ws.Range("C50:D50").VerticalAlignment = -4108
$ws.Range("C50:D50").WrapText = $false
$ws.Range("C50:D50").Locked = $false

# Resize the Excel table / autofilter to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:E159"))
